# feat: add 2022-Q1 data
#
# Starting point: 4 sheets -> 2021-Q1, 2021-Q2, 2021-Q3, 总计.
# Target: 5 sheets -> 2021-Q1, 2021-Q2, 2021-Q3, 2022-Q1, 总计, where:
#   - "2022-Q1" is a new holdings-detail sheet (same shape as the other
#     2021-Qx sheets: 基金代码/基金名称/基金规模/股票总仓位/仓位占比/
#     持有市值(亿元)/仓位排名, one row per fund).
#   - "总计" gains a new leading row for 2022-Q1 (3 holdings, 22.01 亿元),
#     with the previously-existing rows pushed down.
#
# Strategy: duplicate the existing "总计" sheet (so both new sheets start
# out with its styling - bold/bordered/centered header row + styled index
# column - already applied), then rename the two copies and rewrite their
# contents in place. Values are cleared with ClearContents (not Clear) so
# the per-cell formatting survives; formatting for newly-added columns is
# extended by copy/paste-special-formats from an already-styled neighbour
# cell.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")

# Duplicate "总计" -> "总计 (2)" right after itself, giving tab order
# [... , 总计, 总计 (2)]; the original keeps the lower internal sheetId.
$total.Copy($null, $total)

$holdings = $total
$totals = $wb.Worksheets.Item("总计 (2)")

$holdings.Name = "2022-Q1"
$totals.Name = "总计"

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) "2022-Q1" holdings sheet
# ---------------------------------------------------------------------

$holdings.Cells.ClearContents()

# Extend the existing header style (currently only on B1:D1) across to H1.
$holdings.Range("D1").Copy() | Out-Null
$holdings.Range("E1:H1").PasteSpecial($xlPasteFormats)

$holdings.Range("B1").Value = "基金代码"
$holdings.Range("C1").Value = "基金名称"
$holdings.Range("D1").Value = "基金规模"
$holdings.Range("E1").Value = "股票总仓位"
$holdings.Range("F1").Value = "仓位占比"
$holdings.Range("G1").Value = "持有市值(亿元)"
$holdings.Range("H1").Value = "仓位排名"

$holdingsRows = @(
    @{ idx = 0; code = "513050"; name = "易方达中证海外中国互联网50 QDII-ETF"; size = "350.10"; pos = "98.05"; pct = "5.47"; mv = "19.1505"; rank = 4 },
    @{ idx = 1; code = "159605"; name = "广发中证海外中国互联网30（QDII-ETF）"; size = "29.04";  pos = "98.61"; pct = "8.36"; mv = "2.4277";  rank = 4 },
    @{ idx = 2; code = "159607"; name = "嘉实中证海外中国互联网30ETF（QDII）"; size = "5.79";   pos = "98.25"; pct = "7.51"; mv = "0.4348";  rank = 4 }
)

foreach ($row in $holdingsRows) {
    $r = $row.idx + 2
    # Column A (index) already carries the styled xf from the source sheet.
    $holdings.Cells.Item($r, 1).Value = $row.idx
    $holdings.Cells.Item($r, 2).Value = $row.code
    $holdings.Cells.Item($r, 3).Value = $row.name
    # D/E/F/G are stored as plain text in this workbook's convention
    # (e.g. "350.10", "98.05"), not numbers - force text with a leading
    # apostrophe so they don't get auto-coerced to Number.
    $holdings.Cells.Item($r, 4).Value = "'" + $row.size
    $holdings.Cells.Item($r, 5).Value = "'" + $row.pos
    $holdings.Cells.Item($r, 6).Value = "'" + $row.pct
    $holdings.Cells.Item($r, 7).Value = "'" + $row.mv
    $holdings.Cells.Item($r, 8).Value = $row.rank
}

# ---------------------------------------------------------------------
# 2) "总计" totals sheet - prepend the 2022-Q1 row
# ---------------------------------------------------------------------

$totals.Cells.ClearContents()

# A new 5th row is needed (index column only went down to row 4 before);
# extend the styled index-column formatting down one more row.
$totals.Range("A4").Copy() | Out-Null
$totals.Range("A5").PasteSpecial($xlPasteFormats)

$totals.Range("B1").Value = "日期"
$totals.Range("C1").Value = "持有数量(只)"
$totals.Range("D1").Value = "持有市值(亿元)"

$totalsRows = @(
    @{ idx = 0; period = "2022-Q1"; count = 3; mv = 22.01 },
    @{ idx = 1; period = "2021-Q3"; count = 2; mv = 20.08 },
    @{ idx = 2; period = "2021-Q2"; count = 2; mv = 10.67 },
    @{ idx = 3; period = "2021-Q1"; count = 4; mv = 6.16 }
)

foreach ($row in $totalsRows) {
    $r = $row.idx + 2
    $totals.Cells.Item($r, 1).Value = $row.idx
    $totals.Cells.Item($r, 2).Value = $row.period
    $totals.Cells.Item($r, 3).Value = $row.count
    $totals.Cells.Item($r, 4).Value = $row.mv
}

$wb.Worksheets.Item("2021-Q1").Activate()
